# Venture Order List - "View schedule" update
# - renumber / rename the existing work orders (rows 3-6)
# - add three new work orders (rows 7-9)
# - move the active cell selection to D13 (as left by the author when they
#   finished editing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update orderId (column A) on the existing rows 3-6 ---
$ws.Range("A3").Value = 21
$ws.Range("A4").Value = 22
$ws.Range("A5").Value = 23
$ws.Range("A6").Value = 24

# --- Clear the old work-order names so their shared-string slots are freed
#     up (the replacement names get written back further down) ---
$ws.Range("C3:C6").ClearContents() | Out-Null

# --- Copy the date-column formatting (style used by D6:E6) down onto the new
#     rows before the values are written, so the new cells keep the same
#     number format as the rest of the table. ---
$ws.Range("D6:E6").Copy() | Out-Null
$ws.Range("D7:E9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New date values for rows 7-9 ---
$ws.Range("D7").Value = "6/2/2020  12:00PM"
$ws.Range("E7").Value = "06/06/2020 12:00PM"
$ws.Range("D8").Value = "6/2/2020  12:00PM"
$ws.Range("E8").Value = "10/07/2020 12:00PM"
$ws.Range("D9").Value = "6/2/2020  12:00PM"
$ws.Range("E9").Value = "08/08/2020 12:00PM"

# --- Work-order names (rows 3-9) ---
$ws.Range("C3").Value = "Work order 21"
$ws.Range("C4").Value = "Work order 22"
$ws.Range("C5").Value = "Work order 23"
$ws.Range("C6").Value = "Work order 24"
$ws.Range("C7").Value = "work order 25"
$ws.Range("C8").Value = "work order 26"
$ws.Range("C9").Value = "work order 27"

# --- Remaining data for the new rows 7-9 (orderId, lineId, quantity) ---
$ws.Range("A7").Value = 25
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 9520

$ws.Range("A8").Value = 26
$ws.Range("B8").Value = 1
$ws.Range("F8").Value = 5460

$ws.Range("A9").Value = 27
$ws.Range("B9").Value = 1
$ws.Range("F9").Value = 7800

# --- Sheet view / selection changes ---
$ws.Range("D13").Select() | Out-Null
